$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.096.00'
$ws.Range("E2").Value = '  -1.97%  '
$ws.Range("D3").Value = '1.799.96'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.03%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '316.89'
$r.Style = "Normal"
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("E6").Value = '  -0.04%  '
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.5437'
$r.Style = "Normal"
$ws.Range("E7").Value = '  +1.28%  '
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.3802'
$r.Style = "Normal"
$ws.Range("E8").Value = '  +0.83%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.07467'
$r.Style = "Normal"
$ws.Range("E9").Value = '  -0.89%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '41.88'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -1.49%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '1.094'
$r.Style = "Normal"
$ws.Range("E11").Value = '  -1.82%  '
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("E14").Value = '  -2.36%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '7.384'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '1.794.26'
$ws.Range("E16").Value = '  +0.04%  '
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '89.25'
$r.Style = "Normal"
$ws.Range("E17").Value = '  -1.08%  '
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '0.00001065'
$r.Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '0.06519'
$r.Style = "Normal"
$ws.Range("E19").Value = '  +1.27%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '17.43'
$r.Style = "Normal"
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Value = '28.125.09'
$ws.Range("E23").Value = '  -1.84%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '11.17'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -0.68%  '
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '156.88'
$r.Style = "Normal"
$ws.Range("E26").Value = '  -2.34%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '20.39'
$r.Style = "Normal"
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = '2.004.69'
$ws.Range("E28").Value = '  +0.13%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '2.342'
$r.Style = "Normal"
$ws.Range("E29").Value = '  -1.46%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '122.00'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -0.63%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '0.1106'
$r.Style = "Normal"
$ws.Range("E31").Value = '  +7.62%  '
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '1.111'
$r.Style = "Normal"
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("E33").Value = '  -0.54%  '
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '5.546'
$r.Style = "Normal"
$ws.Range("E34").Value = '  -1.87%  '
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '0.06921'
$r.Style = "Normal"
$ws.Range("E35").Value = '  +6.55%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '0.2209'
$r.Style = "Normal"
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("E37").Value = '  -0.57%  '
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '5.082'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +1.01%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '8.425'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -5.51%  '
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '11.18'
$r.Style = "Normal"
$ws.Range("E40").Value = '  -1.33%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.6155'
$r.Style = "Normal"
$ws.Range("E41").Value = '  -1.42%  '
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '1.172'
$r.Style = "Normal"
$ws.Range("E42").Value = '  -3.37%  '
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '1.419'
$r.Style = "Normal"
$ws.Range("E43").Value = '  +1.74%  '
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '13.35'
$r.Style = "Normal"
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("E45").Value = '  +0.56%  '
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.5737'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -2.30%  '
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '124.41'
$r.Style = "Normal"
$ws.Range("E47").Value = '  -1.84%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '1.185'
$r.Style = "Normal"
$ws.Range("E48").Value = '  +2.44%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '1.917'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("E51").Value = '  +38.16%  '
